$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "trajectory"
$ws.Name = "trajectory"

# Add header labels to the first row (columns A, B, C), which were
# previously blank cells.
$ws.Range("A1").Value = "Trajectory"
$ws.Range("B1").Value = "seqId"
$ws.Range("C1").Value = "position"

# Update the active selection to C1 (was I8)
$ws.Range("C1").Select() | Out-Null
